{"js": "async (context) => {\n  // Map of exact \"old text\" -> \"new text\" for every text run that changes.\n  // Every value is unique in the document, so a case-sensitive whole-body\n  // search safely targets the single matching run.\n  const replacements = [\n    [\"2024-06-05 Wednesday\", \"2024-06-06 Thursday\"],\n    [\"12\u00f79=1, 3\", \"22\u00f76=3, 4\"],\n    [\"41\u00f75=8, 1\", \"70\u00f75=14, 0\"],\n    [\"61\u00f73=20, 1\", \"65\u00f75=13, 0\"],\n    [\"71\u00f74=17, 3\", \"63\u00f77=9, 0\"],\n    [\"36\u00f78=4, 4\", \"46\u00f73=15, 1\"],\n    [\"80\u00f73=26, 2\", \"51\u00f78=6, 3\"],\n    [\"87\u00f72=43, 1\", \"60\u00f76=10, 0\"],\n    [\"19\u00f72=9, 1\", \"61\u00f76=10, 1\"],\n    [\"60\u00f77=8, 4\", \"84\u00f72=42, 0\"],\n    [\"10\u00f78=1, 2\", \"66\u00f72=33, 0\"],\n    [\"42\u00f72=21, 0\", \"15\u00f77=2, 1\"],\n    [\"86\u00f73=28, 2\", \"75\u00f75=15, 0\"],\n    [\"92\u00f73=30, 2\", \"36\u00f76=6, 0\"],\n    [\"84\u00f74=21, 0\", \"45\u00f77=6, 3\"],\n    [\"70\u00f77=10, 0\", \"64\u00f74=16, 0\"],\n    [\"60\u00f76=10, 0\", \"89\u00f74=22, 1\"],\n    [\"11\u00f74=2, 3\", \"61\u00f77=8, 5\"],\n    [\"15\u00f79=1, 6\", \"93\u00f73=31, 0\"],\n    [\"18\u00f72=9, 0\", \"36\u00f76=6, 0\"],\n    [\"98\u00f77=14, 0\", \"32\u00f75=6, 2\"],\n    [\"88\u00f79=9, 7\", \"80\u00f74=20, 0\"],\n    [\"23\u00f75=4, 3\", \"51\u00f74=12, 3\"],\n    [\"60\u00f72=30, 0\", \"58\u00f77=8, 2\"],\n    [\"57\u00f75=11, 2\", \"72\u00f74=18, 0\"],\n    [\"69\u00f78=8, 5\", \"92\u00f74=23, 0\"],\n  ];\n\n  const body = context.document.body;\n\n  // First, locate every range to be changed while the document still holds\n  // its ORIGINAL text. Some \"new\" strings coincide with other entries'\n  // \"old\" strings (e.g. \"60\u00f76=10, 0\" and \"36\u00f76=6, 0\" each appear twice\n  // across old/new pairs), so resolving all search hits up-front (rather\n  // than interleaving search+replace) avoids a later search accidentally\n  // matching text that an earlier replacement just produced.\n  const targets = [];\n  for (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length === 0) {\n      throw new Error(`Text not found: \"${oldText}\"`);\n    }\n\n    // Each \"old\" value is unique in the document, so the first (only) hit\n    // is the correct one.\n    targets.push({ range: results.items[0], newText });\n  }\n\n  // Now apply all the replacements.\n  for (const { range, newText } of targets) {\n    range.insertText(newText, \"Replace\");\n  }\n\n  await context.sync();\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# Title paragraph: date line above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-06 Thursday\"\n\n# The worksheet is one 20-row x 5-column table where only every 4th row\n# (1, 5, 9, 13, 17 in 1-based Word indexing) actually holds a division\n# problem; the rows in between are blank spacer rows. Addressing each\n# answer cell directly by (row, column) is unambiguous and sidesteps any\n# issue with duplicate \"before\"/\"after\" text values appearing elsewhere\n# in the table (e.g. \"60\u00f76=10, 0\" is both an old value in one cell and a\n# new value written into a different cell).\n$t = $d.Tables.Item(1)\n\n$answers = @(\n    @{ Row = 1;  Col = 1; Text = \"22\u00f76=3, 4\" },\n    @{ Row = 1;  Col = 2; Text = \"70\u00f75=14, 0\" },\n    @{ Row = 1;  Col = 3; Text = \"65\u00f75=13, 0\" },\n    @{ Row = 1;  Col = 4; Text = \"63\u00f77=9, 0\" },\n    @{ Row = 1;  Col = 5; Text = \"46\u00f73=15, 1\" },\n\n    @{ Row = 5;  Col = 1; Text = \"51\u00f78=6, 3\" },\n    @{ Row = 5;  Col = 2; Text = \"60\u00f76=10, 0\" },\n    @{ Row = 5;  Col = 3; Text = \"61\u00f76=10, 1\" },\n    @{ Row = 5;  Col = 4; Text = \"84\u00f72=42, 0\" },\n    @{ Row = 5;  Col = 5; Text = \"66\u00f72=33, 0\" },\n\n    @{ Row = 9;  Col = 1; Text = \"15\u00f77=2, 1\" },\n    @{ Row = 9;  Col = 2; Text = \"75\u00f75=15, 0\" },\n    @{ Row = 9;  Col = 3; Text = \"36\u00f76=6, 0\" },\n    @{ Row = 9;  Col = 4; Text = \"45\u00f77=6, 3\" },\n    @{ Row = 9;  Col = 5; Text = \"64\u00f74=16, 0\" },\n\n    @{ Row = 13; Col = 1; Text = \"89\u00f74=22, 1\" },\n    @{ Row = 13; Col = 2; Text = \"61\u00f77=8, 5\" },\n    @{ Row = 13; Col = 3; Text = \"93\u00f73=31, 0\" },\n    @{ Row = 13; Col = 4; Text = \"36\u00f76=6, 0\" },\n    @{ Row = 13; Col = 5; Text = \"32\u00f75=6, 2\" },\n\n    @{ Row = 17; Col = 1; Text = \"80\u00f74=20, 0\" },\n    @{ Row = 17; Col = 2; Text = \"51\u00f74=12, 3\" },\n    @{ Row = 17; Col = 3; Text = \"58\u00f77=8, 2\" },\n    @{ Row = 17; Col = 4; Text = \"72\u00f74=18, 0\" },\n    @{ Row = 17; Col = 5; Text = \"92\u00f74=23, 0\" }\n)\n\nforeach ($a in $answers) {\n    $t.Cell($a.Row, $a.Col).Range.Text = $a.Text\n}\n"}
